{"js": "// Locate the table (\u4e2d\u65ad\u4f7f\u7528\u60c5\u51b5 / interrupt usage table) and, within it,\n// the row describing \"SysTick_Handler\", whose last cell (\"\u8bf4\u660e\" / description\n// column) is currently empty. Fill that empty cell's paragraph with two runs:\n// \"1ms\" and \" tick\" (both sized 13pt / w:sz 26, east-asia hinted), producing\n// the final text \"1ms tick\" \u2014 matching the authored OOXML diff exactly.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No tables found in document.\");\n}\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// Find the row whose cells mention \"SysTick_Handler\".\nlet targetRowIndex = -1;\nfor (let r = 0; r < table.values.length; r++) {\n  if (table.values[r].some((cellText) => cellText.indexOf(\"SysTick_Handler\") !== -1)) {\n    targetRowIndex = r;\n    break;\n  }\n}\nif (targetRowIndex === -1) {\n  throw new Error(\"Could not find the SysTick_Handler row.\");\n}\n\n// Target the last column of that row (the \"\u8bf4\u660e\" / description cell).\nconst lastColIndex = table.values[targetRowIndex].length - 1;\nconst cell = table.getCell(targetRowIndex, lastColIndex);\nconst cellBody = cell.body;\ncellBody.load(\"paragraphs\");\nawait context.sync();\n\n// The description cell holds a single (currently empty) paragraph \u2014 this is\n// the paragraph the diff adds the two new runs into.\nconst para = cellBody.paragraphs.items[cellBody.paragraphs.items.length - 1];\n\n// Replace the empty paragraph with an identical one (same pPr) that now also\n// carries the two new runs, matching the diff:\n//   <w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/><w:sz w:val=\"26\"/></w:rPr><w:t>1ms</w:t></w:r>\n//   <w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/><w:sz w:val=\"26\"/></w:rPr><w:t xml:space=\"preserve\"> tick</w:t></w:r>\nconst ooxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p w:rsidR=\"00E96EC1\" w:rsidRPr=\"0098575C\" w:rsidRDefault=\"00E96EC1\" w:rsidP=\"00D31D50\">' +\n  '<w:pPr>' +\n  '<w:spacing w:line=\"220\" w:lineRule=\"atLeast\"/>' +\n  '<w:rPr><w:sz w:val=\"26\"/></w:rPr>' +\n  '</w:pPr>' +\n  '<w:r>' +\n  '<w:rPr><w:rFonts w:hint=\"eastAsia\"/><w:sz w:val=\"26\"/></w:rPr>' +\n  '<w:t>1ms</w:t>' +\n  '</w:r>' +\n  '<w:r>' +\n  '<w:rPr><w:rFonts w:hint=\"eastAsia\"/><w:sz w:val=\"26\"/></w:rPr>' +\n  '<w:t xml:space=\"preserve\"> tick</w:t>' +\n  '</w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\npara.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Locate the table (\u4e2d\u65ad\u4f7f\u7528\u60c5\u51b5 / interrupt usage table) and, within it,\n# the row describing \"SysTick_Handler\", whose last cell (\"\u8bf4\u660e\" / description\n# column) is currently empty. Fill that empty cell's paragraph with two runs:\n# \"1ms\" and \" tick\" (both sized 13pt / w:sz 26, east-asia hinted), producing\n# the final text \"1ms tick\" \u2014 matching the authored OOXML diff exactly.\n\n$d = $word.ActiveDocument\n\n# Find the row whose cells mention \"SysTick_Handler\".\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute(\"SysTick_Handler\")\nif (-not $found) {\n    throw \"Could not find the SysTick_Handler row.\"\n}\n\n$foundCell = $searchRange.Cells.Item(1)\n$tbl = $foundCell.Tables.Item(1)\n\n# Target the last column of that row (the \"\u8bf4\u660e\" / description cell).\n$targetCell = $tbl.Cell($foundCell.RowIndex, $tbl.Columns.Count)\n$cellRange = $targetCell.Range\n\n# Replace the (currently empty) paragraph's contents with the same pPr plus\n# the two new runs, matching the diff:\n#   <w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/><w:sz w:val=\"26\"/></w:rPr><w:t>1ms</w:t></w:r>\n#   <w:r><w:rPr><w:rFonts w:hint=\"eastAsia\"/><w:sz w:val=\"26\"/></w:rPr><w:t xml:space=\"preserve\"> tick</w:t></w:r>\n$xml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p w:rsidR=\"00E96EC1\" w:rsidRPr=\"0098575C\" w:rsidRDefault=\"00E96EC1\" w:rsidP=\"00D31D50\">\n            <w:pPr>\n              <w:spacing w:line=\"220\" w:lineRule=\"atLeast\"/>\n              <w:rPr>\n                <w:sz w:val=\"26\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:hint=\"eastAsia\"/>\n                <w:sz w:val=\"26\"/>\n              </w:rPr>\n              <w:t>1ms</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:hint=\"eastAsia\"/>\n                <w:sz w:val=\"26\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> tick</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$cellRange.InsertXML($xml)\n"}
